$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (row -> A date serial, B, C, D values)
$newRows = @(
    @{ Row = 234; A = 44308; B = 0; C = 0; D = 0 },
    @{ Row = 235; A = 44309; B = 1; C = 1; D = 33.71544167228591 },
    @{ Row = 236; A = 44310; B = 0; C = 1; D = 33.71544167228591 },
    @{ Row = 237; A = 44311; B = 0; C = 1; D = 33.71544167228591 },
    @{ Row = 238; A = 44312; B = 0; C = 1; D = 33.71544167228591 }
)

foreach ($r in $newRows) {
    $srcRow = $r.Row - 1

    # Copy the date-formatted style from the cell directly above into the
    # new cell (so column A keeps its "YYYY-MM-DD HH:MM:SS" date style),
    # then make sure the border survives the copy.
    $ws.Cells.Item($srcRow, 1).Copy()
    $ws.Cells.Item($r.Row, 1).Insert(-4121)
    $ws.Cells.Item($r.Row, 1).Borders.LineStyle = 1

    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
